$d = $word.ActiveDocument

# --- Edit 1: the run that precedes the SECOND inline picture ("Picture 2")
#     is missing <w:noProof/> in its rPr (the run before the first picture
#     already has it). Set Range.NoProofing on that run's range, which
#     serialises to <w:noProof/> right after <w:rFonts.../> in the rPr.
if ($d.InlineShapes.Count -lt 2) {
    throw "Expected at least 2 inline shapes, found " + $d.InlineShapes.Count
}
$secondPicture = $d.InlineShapes.Item(2)
$secondPicture.Range.NoProofing = 1

# --- Edit 2: insert three new paragraphs (a Chinese-caption banner line,
#     plus two bilingual service-description lines mixing MS Mincho / MS
#     Gothic / SimSun runs) right after the second "HOT LINE DICH VU" /
#     "0567 998 998" paragraph, and before the empty paragraph that
#     closes out the document body.
$paraCount = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs.Item($paraCount - 1)
if ($anchorPara.Range.Text.Trim() -ne "0567 998 998") {
    throw "Anchor paragraph mismatch, expected '0567 998 998' but found: " + $anchorPara.Range.Text
}
$insertRange = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

# The fragment's final element is a deliberately empty, namespace-only
# <w:p/>. InsertXML always folds the *last* paragraph of an inserted
# fragment into whichever paragraph follows the insertion point (keeping
# that paragraph's own pPr/formatting and prefixing its runs with any
# runs from the fragment's last paragraph). Making that last fragment
# paragraph empty means nothing gets prefixed into the document's real
# trailing paragraph; it just leaves behind an extra blank paragraph that
# we delete afterwards, so the pre-existing trailing paragraph (and its
# paraId/rsid metadata) stay completely untouched.
$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:beforeAutospacing="1" w:after="0" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Thêm dòng chữ trung quốc dưới banner</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:beforeAutospacing="1" w:after="0" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="MS Mincho" w:eastAsia="MS Mincho" w:hAnsi="MS Mincho" w:cs="MS Mincho" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>保</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:cs="MS Gothic" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>养</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>车，保养产品，洗车</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> ( Chăm Sóc Xe, Lắp Phụ Kiện, Rửa Xe)</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:beforeAutospacing="1" w:after="0" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="MS Mincho" w:eastAsia="MS Mincho" w:hAnsi="MS Mincho" w:cs="MS Mincho" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>修理中心</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MS Mincho" w:eastAsia="MS Mincho" w:hAnsi="MS Mincho" w:cs="MS Mincho" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>保</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:cs="MS Gothic" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>养汽</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun" w:hint="eastAsia"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t>车</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> (Trung Tâm Sửa Chữa-Bảo Dưỡng Ô Tô)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$beforeCount = $d.Paragraphs.Count
$insertRange.InsertXML($fragment)
$afterCount = $d.Paragraphs.Count

# Drop the leftover blank placeholder paragraph created by the absorbed
# last <w:p/> of the fragment (it sits right before the original trailing
# empty paragraph).
$leftover = $d.Paragraphs.Item($afterCount - 1)
if ($leftover.Range.Text.Trim().Length -eq 0) {
    $leftover.Range.Delete()
}

Write-Host "Inserted paragraphs. Count before/after/final:" $beforeCount $afterCount $d.Paragraphs.Count
